# Update cryptocurrency price/volume figures per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.277.03'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.95%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.158.08'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.34%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.32'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.11%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.52'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.71%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.157.33'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.33%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.85%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.164'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.69%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.23'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.67%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.504'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +6.58%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000260'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +13.28%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.46'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.90%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.673.37'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.75%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.260.61'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.98%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.19'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.53%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.160.72'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.64%  '

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.61%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '512.02'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +6.22%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.89'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +6.33%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.730'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +7.29%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.48'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +12.14%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.89'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.83%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.65'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.52%  '

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.03%  '

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.13%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.78'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +7.13%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.18'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.14%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '28.13'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.88%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.11%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.22%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.66'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +5.87%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.06'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +7.99%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.65'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.68%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.80'

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '478.58'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.20%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0424'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.24%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0857'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.85%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.01'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.50%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.131.89'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.31%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.64'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.23%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.97%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.292'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +10.00%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.45'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +12.71%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.33'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.05%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0571'

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.44%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.32'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +10.56%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '118.75'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.58%  '
